# Update computed LR-pair statistics in the active sheet with newly
# recomputed TPM-derived values (see commit: "update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 12.794431
$ws.Range("H2").Value = 38.38329299999999
$ws.Range("I2").Value = 0.8634711242729795
$ws.Range("J2").Value = 0.8634711242729793
$ws.Range("M2").Value = 0.257284
$ws.Range("N2").Value = 0.771852
$ws.Range("O2").Value = 0.05106290078335718
$ws.Range("P2").Value = 0.05106290078335718
$ws.Range("Q2").Value = 3.291802385403999
$ws.Range("R2").Value = 29.626221468636
$ws.Range("S2").Value = 0.04409134034804503
$ws.Range("T2").Value = 0.04409134034804502

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = 12.794431
$ws.Range("H3").Value = 38.38329299999999
$ws.Range("I3").Value = 0.8634711242729795
$ws.Range("J3").Value = 0.8634711242729793
$ws.Range("O3").Value = 0.8646092045957484
$ws.Range("P3").Value = 0.8646092045957485
$ws.Range("Q3").Value = 55.737582441814
$ws.Range("R3").Value = 501.638241976326
$ws.Range("S3").Value = 0.7465650819490574
$ws.Range("T3").Value = 0.7465650819490574

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value = 12.794431
$ws.Range("H4").Value = 38.38329299999999
$ws.Range("I4").Value = 0.8634711242729795
$ws.Range("J4").Value = 0.8634711242729793
$ws.Range("O4").Value = 0.08432789462089441
$ws.Range("P4").Value = 0.08432789462089442
$ws.Range("Q4").Value = 5.436251376451999
$ws.Range("R4").Value = 48.92626238806799
$ws.Range("S4").Value = 0.07281470197587703
$ws.Range("T4").Value = 0.07281470197587703

# Row 5 (FAPs -> ECs)
$ws.Range("H5").Value = 4.00473
$ws.Range("I5").Value = 0.09009046502366876
$ws.Range("J5").Value = 0.09009046502366874
$ws.Range("M5").Value = 0.257284
$ws.Range("N5").Value = 0.771852
$ws.Range("O5").Value = 0.05106290078335718
$ws.Range("P5").Value = 0.05106290078335718
$ws.Range("Q5").Value = 0.3434509844400001
$ws.Range("R5").Value = 3.09105885996
$ws.Range("S5").Value = 0.004600280477030108
$ws.Range("T5").Value = 0.004600280477030107

# Row 6 (FAPs -> FAPs)
$ws.Range("H6").Value = 4.00473
$ws.Range("I6").Value = 0.09009046502366876
$ws.Range("J6").Value = 0.09009046502366874
$ws.Range("O6").Value = 0.8646092045957484
$ws.Range("P6").Value = 0.8646092045957485
$ws.Range("Q6").Value = 5.815393914540001
$ws.Range("R6").Value = 52.33854523086001
$ws.Range("S6").Value = 0.07789304530577533
$ws.Range("T6").Value = 0.07789304530577533

# Row 7 (FAPs -> MuSCs)
$ws.Range("H7").Value = 4.00473
$ws.Range("I7").Value = 0.09009046502366876
$ws.Range("J7").Value = 0.09009046502366874
$ws.Range("O7").Value = 0.08432789462089441
$ws.Range("P7").Value = 0.08432789462089442
$ws.Range("Q7").Value = 0.56719257972
$ws.Range("R7").Value = 5.10473321748
$ws.Range("S7").Value = 0.007597139240863312
$ws.Range("T7").Value = 0.007597139240863312

# Row 8 (MuSCs -> ECs)
$ws.Range("I8").Value = 0.04643841070335186
$ws.Range("J8").Value = 0.04643841070335185
$ws.Range("M8").Value = 0.257284
$ws.Range("N8").Value = 0.771852
$ws.Range("O8").Value = 0.05106290078335718
$ws.Range("P8").Value = 0.05106290078335718
$ws.Range("Q8").Value = 0.1770366915933334
$ws.Range("R8").Value = 1.59333022434
$ws.Range("S8").Value = 0.002371279958282048
$ws.Range("T8").Value = 0.002371279958282047

# Row 9 (MuSCs -> FAPs)
$ws.Range("I9").Value = 0.04643841070335186
$ws.Range("J9").Value = 0.04643841070335185
$ws.Range("O9").Value = 0.8646092045957484
$ws.Range("P9").Value = 0.8646092045957485
$ws.Range("S9").Value = 0.04015107734091574
$ws.Range("T9").Value = 0.04015107734091573

# Row 10 (MuSCs -> MuSCs)
$ws.Range("I10").Value = 0.04643841070335186
$ws.Range("J10").Value = 0.04643841070335185
$ws.Range("O10").Value = 0.08432789462089441
$ws.Range("P10").Value = 0.08432789462089442
$ws.Range("S10").Value = 0.003916053404154071
$ws.Range("T10").Value = 0.003916053404154071

$wb.Save()
